$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "30.385.53"
$ws.Range("E2").Value = "  -0.87%  "
Set-TextValue "D3" "1.859.62"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.21%  "
Set-TextValue "D5" "234.56"
$ws.Range("E5").Value = "  -1.74%  "
Set-TextValue "D6" "1.0000"
$ws.Range("E6").Value = "  -0.18%  "
Set-TextValue "D7" "0.4741"
$ws.Range("E7").Value = "  -1.12%  "
Set-TextValue "D8" "0.2745"
$ws.Range("E8").Value = "  -2.79%  "
Set-TextValue "D9" "0.06440"
$ws.Range("E9").Value = "  -1.04%  "
Set-TextValue "D10" "1.869.09"
$ws.Range("E10").Value = "  -0.79%  "
Set-TextValue "D11" "0.07436"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("E12").Value = "  -0.73%  "
Set-TextValue "D13" "5.009"
$ws.Range("E13").Value = "  -1.49%  "
Set-TextValue "D14" "85.61"
$ws.Range("E14").Value = "  -2.55%  "
Set-TextValue "D15" "0.6360"
$ws.Range("E15").Value = "  -4.06%  "
Set-TextValue "D16" "30.347.55"
$ws.Range("E16").Value = "  -0.81%  "
Set-TextValue "D17" "0.9999"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D18" "12.81"
$ws.Range("E18").Value = "  -3.33%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D19" "231.07"
$ws.Range("E19").Value = "  +2.27%  "
Set-TextValue "D20" "0.000007417"
$ws.Range("E20").Value = "  -1.92%  "
Set-TextValue "D21" "2.102.19"
$ws.Range("E21").Value = "  -2.80%  "
Set-TextValue "D22" "0.9997"
$ws.Range("E22").Value = "  -0.20%  "
Set-TextValue "D23" "5.017"
$ws.Range("E23").Value = "  -4.78%  "
Set-TextValue "D24" "6.009"
$ws.Range("E24").Value = "  -2.16%  "
Set-TextValue "D25" "9.297"
$ws.Range("E25").Value = "  +0.22%  "
Set-TextValue "D26" "165.41"
$ws.Range("E26").Value = "  -1.51%  "
Set-TextValue "D27" "17.97"
$ws.Range("E27").Value = "  -2.86%  "
Set-TextValue "D28" "1.899"
$ws.Range("E28").Value = "  -1.60%  "
Set-TextValue "D29" "0.1044"
$ws.Range("E29").Value = "  +7.94%  "
Set-TextValue "D30" "1.392"
$ws.Range("E30").Value = "  -0.93%  "
Set-TextValue "D31" "4.152"
$ws.Range("E31").Value = "  -4.16%  "
Set-TextValue "D32" "3.941"
$ws.Range("E32").Value = "  -1.36%  "
Set-TextValue "D33" "0.04913"
$ws.Range("E33").Value = "  -2.82%  "
Set-TextValue "D34" "1.154"
$ws.Range("E34").Value = "  -5.53%  "
Set-TextValue "D35" "0.7284"
$ws.Range("E35").Value = "  -2.44%  "
Set-TextValue "D36" "0.9992"
$ws.Range("E36").Value = "  -0.43%  "
Set-TextValue "D37" "2.695"
$ws.Range("E37").Value = "  -0.75%  "
Set-TextValue "D38" "0.01898"
$ws.Range("E38").Value = "  +1.93%  "
Set-TextValue "D39" "2.654"
$ws.Range("E39").Value = "  +0.75%  "
Set-TextValue "D40" "0.9115"
$ws.Range("E40").Value = "  -0.27%  "
Set-TextValue "D41" "1.975"
$ws.Range("E41").Value = "  -4.50%  "
Set-TextValue "D42" "105.73"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E43").Value = "  +0.07%  "
Set-TextValue "D44" "0.4123"
$ws.Range("E44").Value = "  -3.23%  "
Set-TextValue "D45" "5.573"
$ws.Range("E45").Value = "  -3.06%  "
Set-TextValue "D46" "7.151"
Set-TextValue "D47" "61.16"
$ws.Range("E47").Value = "  -4.63%  "
Set-TextValue "D49" "8.724"
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D50" "33.45"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D51" "1.410"
$ws.Range("E51").Value = "  -3.98%  "
